# Add one more neural-network test run (test #13) to the log table,
# then extend the defined name and the "Accuracy" line-chart series so
# they include the new row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- 1. Append the new data row (row 14) ------------------------------
$newRow = 14
$ws.Cells.Item($newRow, 1).Value = 13      # TestNo
$ws.Cells.Item($newRow, 2).Value = 1000    # TrainAmount
$ws.Cells.Item($newRow, 3).Value = 3       # HiddenLayerNodes
$ws.Cells.Item($newRow, 4).Value = 0.2     # LearningRate
$ws.Cells.Item($newRow, 5).Value = 100     # Epochs
$ws.Cells.Item($newRow, 6).Value = 100     # TestAmount
$ws.Cells.Item($newRow, 7).Value = 96      # Accuracy

# --- 2. Extend the "neuralNetworkLog" defined name to the new range ---
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Sheet1!neuralNetworkLog") {
        $n.RefersTo = "=Sheet1!`$A`$1:`$G`$14"
    }
}

# --- 3. Extend the chart series so it plots the new Accuracy value ----
$co = $ws.ChartObjects(1)
$chart = $co.Chart
$series = $chart.SeriesCollection(1)
$series.Formula = "=SERIES(Sheet1!`$G`$1,,Sheet1!`$G`$2:`$G`$14,1)"
